$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 8.540560666666666
$ws.Range("H2").Value = 25.621682
$ws.Range("I2").Value = 0.4159358086620884
$ws.Range("J2").Value = 0.4159358086620884
$ws.Range("M2").Value = 9.031965666666666
$ws.Range("N2").Value = 27.095897
$ws.Range("O2").Value = 0.4424406034784756
$ws.Range("P2").Value = 0.4424406034784755
$ws.Range("Q2").Value = 77.13805071541711
$ws.Range("R2").Value = 694.2424564387539
$ws.Range("S2").Value = 0.1840268901927621
$ws.Range("T2").Value = 0.1840268901927621
$ws.Range("G3").Value = 8.540560666666666
$ws.Range("H3").Value = 25.621682
$ws.Range("I3").Value = 0.4159358086620884
$ws.Range("J3").Value = 0.4159358086620884
$ws.Range("O3").Value = 0.4469933372071527
$ws.Range("P3").Value = 0.4469933372071526
$ws.Range("Q3").Value = 77.93180473006998
$ws.Range("R3").Value = 701.3862425706299
$ws.Range("S3").Value = 0.1859205351778226
$ws.Range("T3").Value = 0.1859205351778226
$ws.Range("G4").Value = 8.540560666666666
$ws.Range("H4").Value = 25.621682
$ws.Range("I4").Value = 0.4159358086620884
$ws.Range("J4").Value = 0.4159358086620884
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.2635683333333333
$ws.Range("N4").Value = 0.790705
$ws.Range("O4").Value = 0.01291117977653399
$ws.Range("P4").Value = 0.01291117977653399
$ws.Range("Q4").Value = 2.251021340645555
$ws.Range("R4").Value = 20.25919206581
$ws.Range("S4").Value = 0.005370222001134268
$ws.Range("T4").Value = 0.005370222001134267
$ws.Range("G5").Value = 8.540560666666666
$ws.Range("H5").Value = 25.621682
$ws.Range("I5").Value = 0.4159358086620884
$ws.Range("J5").Value = 0.4159358086620884
$ws.Range("M5").Value = 1.809602666666667
$ws.Range("N5").Value = 5.428808
$ws.Range("O5").Value = 0.08864534315615299
$ws.Range("P5").Value = 0.08864534315615297
$ws.Range("Q5").Value = 15.45502135722844
$ws.Range("R5").Value = 139.095192215056
$ws.Range("S5").Value = 0.03687077248978282
$ws.Range("T5").Value = 0.03687077248978281
$ws.Range("G6").Value = 8.540560666666666
$ws.Range("H6").Value = 25.621682
$ws.Range("I6").Value = 0.4159358086620884
$ws.Range("J6").Value = 0.4159358086620884
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.1839203333333334
$ws.Range("N6").Value = 0.5517610000000001
$ws.Range("O6").Value = 0.009009536381684918
$ws.Range("P6").Value = 0.009009536381684917
$ws.Range("Q6").Value = 1.570782764666889
$ws.Range("R6").Value = 14.137044882002
$ws.Range("S6").Value = 0.003747388800586623
$ws.Range("T6").Value = 0.003747388800586622
$ws.Range("I7").Value = 0.563694901924408
$ws.Range("J7").Value = 0.563694901924408
$ws.Range("M7").Value = 9.031965666666666
$ws.Range("N7").Value = 27.095897
$ws.Range("O7").Value = 0.4424406034784756
$ws.Range("P7").Value = 0.4424406034784755
$ws.Range("Q7").Value = 104.5409532603928
$ws.Range("R7").Value = 940.868579343535
$ws.Range("S7").Value = 0.2494015125851752
$ws.Range("T7").Value = 0.2494015125851751
$ws.Range("I8").Value = 0.563694901924408
$ws.Range("J8").Value = 0.563694901924408
$ws.Range("O8").Value = 0.4469933372071527
$ws.Range("P8").Value = 0.4469933372071526
$ws.Range("S8").Value = 0.2519678653778498
$ws.Range("T8").Value = 0.2519678653778497
$ws.Range("I9").Value = 0.563694901924408
$ws.Range("J9").Value = 0.563694901924408
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.2635683333333333
$ws.Range("N9").Value = 0.790705
$ws.Range("O9").Value = 0.01291117977653399
$ws.Range("P9").Value = 0.01291117977653399
$ws.Range("Q9").Value = 3.050685291863889
$ws.Range("R9").Value = 27.456167626775
$ws.Range("S9").Value = 0.007277966217861729
$ws.Range("T9").Value = 0.007277966217861727
$ws.Range("I10").Value = 0.563694901924408
$ws.Range("J10").Value = 0.563694901924408
$ws.Range("M10").Value = 1.809602666666667
$ws.Range("N10").Value = 5.428808
$ws.Range("O10").Value = 0.08864534315615299
$ws.Range("P10").Value = 0.08864534315615297
$ws.Range("Q10").Value = 20.94533956147111
$ws.Range("R10").Value = 188.50805605324
$ws.Range("S10").Value = 0.04996892801646315
$ws.Range("T10").Value = 0.04996892801646313
$ws.Range("I11").Value = 0.563694901924408
$ws.Range("J11").Value = 0.563694901924408
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 0.6666666666666666
$ws.Range("M11").Value = 0.1839203333333334
$ws.Range("N11").Value = 0.5517610000000001
$ws.Range("O11").Value = 0.009009536381684918
$ws.Range("P11").Value = 0.009009536381684917
$ws.Range("Q11").Value = 2.128795400717222
$ws.Range("R11").Value = 19.159158606455
$ws.Range("S11").Value = 0.005078629727058265
$ws.Range("T11").Value = 0.005078629727058265
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = 0.3333333333333333
$ws.Range("G12").Value = 0.026642
$ws.Range("H12").Value = 0.079926
$ws.Range("I12").Value = 0.001297498167494471
$ws.Range("J12").Value = 0.001297498167494471
$ws.Range("M12").Value = 9.031965666666666
$ws.Range("N12").Value = 27.095897
$ws.Range("O12").Value = 0.4424406034784756
$ws.Range("P12").Value = 0.4424406034784755
$ws.Range("Q12").Value = 0.2406296292913333
$ws.Range("R12").Value = 2.165666663622
$ws.Range("S12").Value = 0.00057406587223847
$ws.Range("T12").Value = 0.0005740658722384699
$ws.Range("E13").Value = 1
$ws.Range("F13").Value = 0.3333333333333333
$ws.Range("G13").Value = 0.026642
$ws.Range("H13").Value = 0.079926
$ws.Range("I13").Value = 0.001297498167494471
$ws.Range("J13").Value = 0.001297498167494471
$ws.Range("O13").Value = 0.4469933372071527
$ws.Range("P13").Value = 0.4469933372071526
$ws.Range("Q13").Value = 0.2431057190099999
$ws.Range("R13").Value = 2.187951471089999
$ws.Range("S13").Value = 0.0005799730359085188
$ws.Range("T13").Value = 0.0005799730359085188
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 0.3333333333333333
$ws.Range("G14").Value = 0.026642
$ws.Range("H14").Value = 0.079926
$ws.Range("I14").Value = 0.001297498167494471
$ws.Range("J14").Value = 0.001297498167494471
$ws.Range("K14").Value = 2
$ws.Range("L14").Value = 0.6666666666666666
$ws.Range("M14").Value = 0.2635683333333333
$ws.Range("N14").Value = 0.790705
$ws.Range("O14").Value = 0.01291117977653399
$ws.Range("P14").Value = 0.01291117977653399
$ws.Range("Q14").Value = 0.007021987536666667
$ws.Range("R14").Value = 0.06319788782999999
$ws.Range("S14").Value = 0.00001675223210024453
$ws.Range("T14").Value = 0.00001675223210024453
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 0.3333333333333333
$ws.Range("G15").Value = 0.026642
$ws.Range("H15").Value = 0.079926
$ws.Range("I15").Value = 0.001297498167494471
$ws.Range("J15").Value = 0.001297498167494471
$ws.Range("M15").Value = 1.809602666666667
$ws.Range("N15").Value = 5.428808
$ws.Range("O15").Value = 0.08864534315615299
$ws.Range("P15").Value = 0.08864534315615297
$ws.Range("Q15").Value = 0.04821143424533333
$ws.Range("R15").Value = 0.433902908208
$ws.Range("S15").Value = 0.0001150171703020271
$ws.Range("T15").Value = 0.000115017170302027
$ws.Range("E16").Value = 1
$ws.Range("F16").Value = 0.3333333333333333
$ws.Range("G16").Value = 0.026642
$ws.Range("H16").Value = 0.079926
$ws.Range("I16").Value = 0.001297498167494471
$ws.Range("J16").Value = 0.001297498167494471
$ws.Range("K16").Value = 2
$ws.Range("L16").Value = 0.6666666666666666
$ws.Range("M16").Value = 0.1839203333333334
$ws.Range("N16").Value = 0.5517610000000001
$ws.Range("O16").Value = 0.009009536381684918
$ws.Range("P16").Value = 0.009009536381684917
$ws.Range("Q16").Value = 0.004900005520666667
$ws.Range("R16").Value = 0.044100049686
$ws.Range("S16").Value = 0.00001168985694521095
$ws.Range("T16").Value = 0.00001168985694521095
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 0.391608
$ws.Range("H17").Value = 1.174824
$ws.Range("I17").Value = 0.01907179124600912
$ws.Range("J17").Value = 0.01907179124600912
$ws.Range("M17").Value = 9.031965666666666
$ws.Range("N17").Value = 27.095897
$ws.Range("O17").Value = 0.4424406034784756
$ws.Range("P17").Value = 0.4424406034784755
$ws.Range("Q17").Value = 3.536990010792
$ws.Range("R17").Value = 31.832910097128
$ws.Range("S17").Value = 0.008438134828299781
$ws.Range("T17").Value = 0.008438134828299781
$ws.Range("E18").Value = 3
$ws.Range("F18").Value = 1
$ws.Range("G18").Value = 0.391608
$ws.Range("H18").Value = 1.174824
$ws.Range("I18").Value = 0.01907179124600912
$ws.Range("J18").Value = 0.01907179124600912
$ws.Range("O18").Value = 0.4469933372071527
$ws.Range("P18").Value = 0.4469933372071526
$ws.Range("Q18").Value = 3.573385797239999
$ws.Range("R18").Value = 32.16047217516
$ws.Range("S18").Value = 0.008524963615571776
$ws.Range("T18").Value = 0.008524963615571774
$ws.Range("E19").Value = 3
$ws.Range("F19").Value = 1
$ws.Range("G19").Value = 0.391608
$ws.Range("H19").Value = 1.174824
$ws.Range("I19").Value = 0.01907179124600912
$ws.Range("J19").Value = 0.01907179124600912
$ws.Range("K19").Value = 2
$ws.Range("L19").Value = 0.6666666666666666
$ws.Range("M19").Value = 0.2635683333333333
$ws.Range("N19").Value = 0.790705
$ws.Range("O19").Value = 0.01291117977653399
$ws.Range("P19").Value = 0.01291117977653399
$ws.Range("Q19").Value = 0.10321546788
$ws.Range("R19").Value = 0.9289392109200001
$ws.Range("S19").Value = 0.000246239325437751
$ws.Range("T19").Value = 0.0002462393254377509
$ws.Range("E20").Value = 3
$ws.Range("F20").Value = 1
$ws.Range("G20").Value = 0.391608
$ws.Range("H20").Value = 1.174824
$ws.Range("I20").Value = 0.01907179124600912
$ws.Range("J20").Value = 0.01907179124600912
$ws.Range("M20").Value = 1.809602666666667
$ws.Range("N20").Value = 5.428808
$ws.Range("O20").Value = 0.08864534315615299
$ws.Range("P20").Value = 0.08864534315615297
$ws.Range("Q20").Value = 0.708654881088
$ws.Range("R20").Value = 6.377893929792001
$ws.Range("S20").Value = 0.001690625479604993
$ws.Range("T20").Value = 0.001690625479604992
$ws.Range("E21").Value = 3
$ws.Range("F21").Value = 1
$ws.Range("G21").Value = 0.391608
$ws.Range("H21").Value = 1.174824
$ws.Range("I21").Value = 0.01907179124600912
$ws.Range("J21").Value = 0.01907179124600912
$ws.Range("K21").Value = 2
$ws.Range("L21").Value = 0.6666666666666666
$ws.Range("M21").Value = 0.1839203333333334
$ws.Range("N21").Value = 0.5517610000000001
$ws.Range("O21").Value = 0.009009536381684918
$ws.Range("P21").Value = 0.009009536381684917
$ws.Range("Q21").Value = 0.07202467389600001
$ws.Range("R21").Value = 0.6482220650640002
$ws.Range("S21").Value = 0.0001718279970948191
$ws.Range("T21").Value = 0.0001718279970948191
